$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Otago row (row 4) totals
$ws.Range("B4").Value = 14
$ws.Range("C4").Value = 13

# Rename current row 14 "Southland" to "South Canterbury"
$ws.Range("A14").Value = "South Canterbury"

# Insert a new row before row 15 (shifts old row 15 "Tasman" down to row 16)
$ws.Rows.Item(15).Insert()

# Populate the newly inserted row 15 with "Southland" data
$ws.Range("A15").Value = "Southland"
$ws.Range("B15").Value = 1
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 0
